$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date strings from slash format (DD/MM/YYYY) to dash format (DD-MM-YYYY)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
}

# Update the numeric attendance values (D=4, E=5, G=7, H=8) that changed
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 8).Value = 0

$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 8).Value = 0

$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 8).Value = 0
